$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B to "No" for every test row except row 2 (LoginSetup) and
# row 12 (DataNexusDataProfiler, already "No"). Rows 3-11 and 13-28 flip
# from "Yes" to "No".
foreach ($r in 3..11) {
    $ws.Cells.Item($r, 2).Value = "No"
}
foreach ($r in 13..28) {
    $ws.Cells.Item($r, 2).Value = "No"
}

# Update the selection / scroll position recorded in the sheet view.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B8").Select()
